$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.608.86"
$ws.Range("E2").Value = "  -2.33%  "

$ws.Range("D3").Value = "1.820.84"
$ws.Range("E3").Value = "  -1.86%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.77%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.01"
$ws.Range("E5").Value = "  -1.55%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4570"
$ws.Range("E7").Value = "  -1.43%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3675"
$ws.Range("E8").Value = "  -1.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07166"
$ws.Range("E9").Value = "  -1.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8803"
$ws.Range("E10").Value = "  -0.78%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07785"
$ws.Range("E11").Value = "  -1.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.44"
$ws.Range("E12").Value = "  -3.36%  "

$ws.Range("D13").Value = "1.777.39"
$ws.Range("E13").Value = "  -1.85%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.304"
$ws.Range("E14").Value = "  -1.72%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.388"
$ws.Range("E15").Value = "  -2.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "86.55"
$ws.Range("E16").Value = "  -4.99%  "

$ws.Range("E17").Value = "  +0.93%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008622"
$ws.Range("E18").Value = "  -3.53%  "

$ws.Range("E19").Value = "  +0.75%  "

$ws.Range("D20").Value = "26.687.32"
$ws.Range("E20").Value = "  -2.13%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.28"
$ws.Range("E21").Value = "  -2.97%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.005"
$ws.Range("E22").Value = "  -1.60%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.49"
$ws.Range("E23").Value = "  -0.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.993"
$ws.Range("E24").Value = "  +2.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.70"
$ws.Range("E25").Value = "  +0.18%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.01"
$ws.Range("E26").Value = "  -2.35%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.080"
$ws.Range("E27").Value = "  +1.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "113.19"
$ws.Range("E28").Value = "  -2.41%  "

$ws.Range("E29").Value = "  -3.63%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08704"
$ws.Range("E30").Value = "  -1.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.069"
$ws.Range("E31").Value = "  -2.25%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.523"
$ws.Range("E32").Value = "  +0.22%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7377"
$ws.Range("E33").Value = "  -4.09%  "

$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.712"
$ws.Range("E34").Value = "  -0.32%  "

$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.125"
$ws.Range("E35").Value = "  -3.70%  "

$ws.Range("E36").Value = "  +0.64%  "

$ws.Range("E37").Value = "  -2.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01948"
$ws.Range("E38").Value = "  +0.42%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05129"
$ws.Range("E39").Value = "  -1.79%  "

$ws.Range("E40").Value = "  -1.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.034"
$ws.Range("E41").Value = "  -0.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5049"
$ws.Range("E42").Value = "  -1.49%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1560"
$ws.Range("E43").Value = "  -4.21%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.197"
$ws.Range("E44").Value = "  -2.85%  "

$ws.Range("E45").Value = "  +0.89%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4640"
$ws.Range("E46").Value = "  -3.27%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.981"
$ws.Range("E47").Value = "  -3.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "101.51"
$ws.Range("E48").Value = "  -1.53%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.597"
$ws.Range("E49").Value = "  -2.94%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06037"
$ws.Range("E50").Value = "  -2.73%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.55"
$ws.Range("E51").Value = "  -1.60%  "
